$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.10955972832100969,
    -0.0059999999480808697,
    -0.0039999999590749624,
    -0.0079999999237028163,
    -0.0029999999650778264,
    -0.0019999999722628559,
    -0.0099999999050637278,
    -0.0099999999058293376,
    -0.0019999999768423038,
    -0.0019999999801321167,
    -0.0029999999718075543,
    -0.0034999999683855698,
    0.01907153387886229,
    -0.0079999999379207765,
    -0.00099999999889277547,
    0.020936736126373923,
    -0.0019999999915647493,
    -0.003999999974310775,
    -0.0028050937951658028,
    -0.0039999999620778937,
    -0.0039999999616293636,
    -0.0039999999612696513,
    -0.0049999999477767787,
    -0.019999999814903191,
    -0.019999999812169378,
    -0.0024999999677657314,
    -0.0024999999673829265,
    -0.00199999997024225,
    -0.0069999999272054581,
    -0.059999999477032606,
    -0.0069999999377525768,
    0.069941757791172066,
    -0.0039999999665134567
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
